# edit.ps1 -- Junction_Flooding_413 dataset refresh (custom accuracy + 1000 new data points)
# Replaces the 4 sample data rows with a new set of readings, widens several
# data columns by one character, and drops the now-unused 6th data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the sample readings in rows 2-5 -------------------------------
# Row 2
$ws.Range("A2").Value = 45094.50694444445
$ws.Range("B2").Value = 14.315
$ws.Range("C2").Value = 9.452
$ws.Range("D2").Value = 3.548
$ws.Range("E2").Value = 30.929
$ws.Range("F2").Value = 23.407
$ws.Range("G2").Value = 11.081
$ws.Range("H2").Value = 33.367
$ws.Range("I2").Value = 17.452
$ws.Range("J2").Value = 7.022
$ws.Range("K2").Value = 10.47
$ws.Range("L2").Value = 12.134
$ws.Range("M2").Value = 12.887
$ws.Range("N2").Value = 3.618
$ws.Range("O2").Value = 11.279
$ws.Range("P2").Value = 15.513
$ws.Range("Q2").Value = 9.993
$ws.Range("R2").Value = 3.077
$ws.Range("S2").Value = 1.709
$ws.Range("T2").Value = 164.3
$ws.Range("U2").Value = 31.24
$ws.Range("V2").Value = 10.411
$ws.Range("W2").Value = 20.071
$ws.Range("X2").Value = 10.298
$ws.Range("Y2").Value = 2.945
$ws.Range("Z2").Value = 17.697
$ws.Range("AA2").Value = 9.196
$ws.Range("AB2").Value = 8.398
$ws.Range("AC2").Value = 9.994999999999999
$ws.Range("AD2").Value = 12.384
$ws.Range("AE2").Value = 3.078
$ws.Range("AF2").Value = 30.175
$ws.Range("AG2").Value = 5.636
$ws.Range("AH2").Value = 13.015

# Row 3
$ws.Range("A3").Value = 45094.51388888889
$ws.Range("B3").Value = 17.21
$ws.Range("C3").Value = 12.379
$ws.Range("D3").Value = 1.787
$ws.Range("E3").Value = 37.629
$ws.Range("F3").Value = 29.906
$ws.Range("G3").Value = 13.423
$ws.Range("H3").Value = 50.882
$ws.Range("I3").Value = 20.942
$ws.Range("J3").Value = 9.055
$ws.Range("K3").Value = 13.294
$ws.Range("L3").Value = 14.994
$ws.Range("M3").Value = 15.939
$ws.Range("N3").Value = 4.347
$ws.Range("O3").Value = 13.535
$ws.Range("P3").Value = 19.062
$ws.Range("Q3").Value = 11.707
$ws.Range("R3").Value = 1.464
$ws.Range("S3").Value = 1.048
$ws.Range("T3").Value = 198.69
$ws.Range("U3").Value = 37.842
$ws.Range("V3").Value = 12.493
$ws.Range("W3").Value = 25.043
$ws.Range("X3").Value = 13.075
$ws.Range("Y3").Value = 2.412
$ws.Range("Z3").Value = 25.224
$ws.Range("AA3").Value = 11.035
$ws.Range("AB3").Value = 9.928000000000001
$ws.Range("AC3").Value = 11.704
$ws.Range("AD3").Value = 15.516
$ws.Range("AE3").Value = 1.199
$ws.Range("AF3").Value = 46.641
$ws.Range("AG3").Value = 6.918
$ws.Range("AH3").Value = 15.619

# Row 4
$ws.Range("A4").Value = 45094.52083333334
$ws.Range("B4").Value = 8.571999999999999
$ws.Range("C4").Value = 6.114
$ws.Range("D4").Value = 1.041
$ws.Range("E4").Value = 18.853
$ws.Range("F4").Value = 14.74
$ws.Range("G4").Value = 6.66
$ws.Range("H4").Value = 29.31
$ws.Range("I4").Value = 10.471
$ws.Range("J4").Value = 4.495
$ws.Range("K4").Value = 6.478
$ws.Range("L4").Value = 7.507
$ws.Range("M4").Value = 8.02
$ws.Range("N4").Value = 2.177
$ws.Range("O4").Value = 6.767
$ws.Range("P4").Value = 9.507999999999999
$ws.Range("Q4").Value = 5.973
$ws.Range("R4").Value = 0.9340000000000001
$ws.Range("S4").Value = 0.582
$ws.Range("T4").Value = 95.696
$ws.Range("U4").Value = 19.055
$ws.Range("V4").Value = 6.247
$ws.Range("W4").Value = 12.509
$ws.Range("X4").Value = 6.502
$ws.Range("Y4").Value = 1.308
$ws.Range("Z4").Value = 13.91
$ws.Range("AA4").Value = 5.517
$ws.Range("AB4").Value = 5.025
$ws.Range("AC4").Value = 5.91
$ws.Range("AD4").Value = 7.728
$ws.Range("AE4").Value = 0.747
$ws.Range("AF4").Value = 27.035
$ws.Range("AG4").Value = 3.405
$ws.Range("AH4").Value = 7.81

# Row 5
$ws.Range("A5").Value = 45094.52777777778
$ws.Range("B5").Value = 2.81
$ws.Range("C5").Value = 1.89
$ws.Range("D5").Value = 0.64
$ws.Range("E5").Value = 6.31
$ws.Range("F5").Value = 4.57
$ws.Range("G5").Value = 2.15
$ws.Range("H5").Value = 12.6
$ws.Range("I5").Value = 3.49
$ws.Range("J5").Value = 1.46
$ws.Range("K5").Value = 1.91
$ws.Range("L5").Value = 2.49
$ws.Range("M5").Value = 2.71
$ws.Range("N5").Value = 0.74
$ws.Range("O5").Value = 2.26
$ws.Range("P5").Value = 3.15
$ws.Range("Q5").Value = 2.16
$ws.Range("R5").Value = 0.67
$ws.Range("S5").Value = 0.32
$ws.Range("T5").Value = 27.05
$ws.Range("U5").Value = 6.55
$ws.Range("V5").Value = 2.08
$ws.Range("W5").Value = 4.19
$ws.Range("X5").Value = 2.12
$ws.Range("Y5").Value = 0.61
$ws.Range("Z5").Value = 5.74
$ws.Range("AA5").Value = 1.84
$ws.Range("AB5").Value = 1.76
$ws.Range("AC5").Value = 2.06
$ws.Range("AD5").Value = 2.51
$ws.Range("AE5").Value = 0.54
$ws.Range("AF5").Value = 11.86
$ws.Range("AG5").Value = 1.06
$ws.Range("AH5").Value = 2.61

# --- Widen columns whose content now needs an extra character of width -----
$ws.Range("B1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("C1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("G1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("I1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("K1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("L1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("M1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("O1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("P1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("Q1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("V1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("W1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("X1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("Z1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AA1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AC1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AD1").EntireColumn.ColumnWidth = 7.166666666666667
$ws.Range("AE1").EntireColumn.ColumnWidth = 6.166666666666667
$ws.Range("AH1").EntireColumn.ColumnWidth = 7.166666666666667

# --- Drop the old 6th data row; used range becomes A1:AH5 ------------------
$ws.Rows(6).Delete()
